# Applies the "Finish accounting for debit/trans/exh" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4: Deadline/estimate increase ---
$ws.Range("H4").Value2 = 150

# --- Row 13: brand-new "DOCS" entry (copy formatting that matches row 12's pattern) ---
$ws.Range("A12:H12").Copy() | Out-Null
$ws.Range("A13:H13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value2 = "DOCS"
$ws.Range("C13").Value2 = "In biên bảng giao hàng"
$ws.Range("D13").Value2 = "Export Excel"
$ws.Range("E13").ClearContents() | Out-Null
$ws.Range("H13").Value2 = 50

# --- Row 14: "Tạm ứng" group head; becomes the top of a new H14:H17 merge ---
$ws.Range("A14").Value2 = 4
$ws.Range("B14").Value2 = "Tạm ứng"
$ws.Range("C14").Value2 = "Phieu de nghi tam ung"
$ws.Range("E14").ClearContents() | Out-Null
$ws.Range("H14").Value2 = 500

# --- Row 15 ---
$ws.Range("C15").Value2 = "Phieu de nghi thanh toan"
$ws.Range("E15").ClearContents() | Out-Null

# --- Row 16 ---
$ws.Range("C16").Value2 = "Tong hop tam ung"
$ws.Range("E16").ClearContents() | Out-Null

# --- Row 17: now holds "Form phieu thu chi" (replacing the removed "customs declaration" row) ---
$ws.Range("A17").ClearContents() | Out-Null
$ws.Range("B17").ClearContents() | Out-Null
$ws.Range("C17").Value2 = "Form phiếu thu chi"
$ws.Range("E17").ClearContents() | Out-Null

# Copy the H9/H10/H11 merge-edge formatting onto the new H14:H17 merge group
# before merging, so borders render the same way as the existing H9:H11 group.
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("H14").Value2 = 500
$ws.Range("H14:H17").Merge() | Out-Null

# --- Row 18: "Others" / "Shipment control" group (copy formatting from old row 16) ---
$ws.Range("A16:B16").Copy() | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A18").Value2 = 5
$ws.Range("B18").Value2 = "Others"
$ws.Range("C18").Value2 = "Shipment control"
$ws.Range("E18").ClearContents() | Out-Null
$ws.Range("H18").Value2 = 300

# --- Row 19: "Kế hoạch vận tải" ---
$ws.Range("C19").Value2 = "Kế hoạch vận tải"
$ws.Range("E19").ClearContents() | Out-Null
$ws.Range("H19").Value2 = 100

# --- View: selection moves to D12, and the saved top-left scroll resets to A1 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select() | Out-Null
